$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "16.80", thousand-dot prices) are preserved exactly as text,
# matching the original inline-string cell contents.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.389.15"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.03"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3759"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.09"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3418"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.165"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07683"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.37"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.992"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.932"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001142"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.574.35"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.35"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06719"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.80"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.246"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5281"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.383.73"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.393"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.776"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "144.52"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.081"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.20"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.749.23"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +8.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.240"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.07"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08525"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02561"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06550"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.294"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.25%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6446"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.07"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6024"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.778"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +11.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.02"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.92%  "
